$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169685244560242
$ws.Range("B1").Value = 2.440001010894775
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.362849712371826
$ws.Range("E1").Value = 1.238023042678833
